# Apply updated Betfair back/lay odds to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("F2").Value = 1.66
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 5.1
$ws.Range("I2").Value = 5.7
$ws.Range("K2").Value = 4.7
$ws.Range("V2").Value = 1.21
$ws.Range("W2").Value = 2.36
$ws.Range("AD2").Value = 21
$ws.Range("AH2").Value = 19

# Row 3
$ws.Range("F3").Value = 1.39
$ws.Range("G3").Value = 1.49
$ws.Range("H3").Value = 9.6
$ws.Range("I3").Value = 13.5
$ws.Range("J3").Value = 4.4
$ws.Range("K3").Value = 5.2
$ws.Range("P3").Value = 1.81
$ws.Range("Q3").Value = 1.99

# Row 4
$ws.Range("G4").Value = 2.12
$ws.Range("I4").Value = 3.85
$ws.Range("J4").Value = 4.1
$ws.Range("K4").Value = 5.6
$ws.Range("Q4").Value = 1.4

# Row 5
$ws.Range("F5").Value = 1.42
$ws.Range("K5").Value = 6.4
$ws.Range("Q5").Value = 1.3

# Row 6
$ws.Range("I6").Value = 3.9
$ws.Range("J6").Value = 3.6

# Row 7
$ws.Range("F7").Value = 1.43
$ws.Range("G7").Value = 1.52
$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 5.1
$ws.Range("K7").Value = 6
$ws.Range("P7").Value = 2.92
$ws.Range("Q7").Value = 1.42

# Row 8
$ws.Range("G8").Value = 2.92
$ws.Range("I8").Value = 2.82
$ws.Range("J8").Value = 3.45
$ws.Range("P8").Value = 1.94
$ws.Range("Q8").Value = 1.9
